$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix up formatting for rows whose role/category changed ---
# (day-04 grew: separator+day-05-header shift from rows 42/43 to rows 44/45)
$ws.Range("A25:B25").Copy() | Out-Null
$ws.Range("A42:B42").PasteSpecial(-4122) | Out-Null
$ws.Range("A25:B25").Copy() | Out-Null
$ws.Range("A43:B43").PasteSpecial(-4122) | Out-Null
$ws.Range("A27:B27").Copy() | Out-Null
$ws.Range("A44:B44").PasteSpecial(-4122) | Out-Null
$ws.Range("A32:B32").Copy() | Out-Null
$ws.Range("A45:B45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 2: set cell values (Folder/Document/Sub-Folder/Description text) ---
$ws.Range("B18").Value = "Contains slides for Angular Unit Testing"
$ws.Range("B20").Value = "Contains important links to various Angular articles and YouTube videos. Participants can use these links as reference"

$ws.Range("A24").Value = "day-01/"
$ws.Range("B24").Value = ""
$ws.Range("A25").Value = "ts-examples"
$ws.Range("B25").Value = "Contains code examples for various concepts in TypeScript. Files have been named appropriately after the concept names. For e.g., 'template-strings.ts' contains demo code for template literal (also called template sting) of TypeScript / ES6"
$ws.Range("A26").Value = "my-app"
$ws.Range("B26").Value = "This app is the first and basic Angular app of the course. This is created to understand:`n1. How to create and run an Angular app using Angular CLI (https://angular.io/cli)`n2. Angular project structure (https://angular.io/guide/file-structure)`n3. Bootstrapping process (https://angular.io/guide/bootstrapping)"
$ws.Range("A27").Value = ""
$ws.Range("B27").Value = ""
$ws.Range("A28").Value = "day-02/"
$ws.Range("B28").Value = ""
$ws.Range("A29").Value = "ts-examples"
$ws.Range("B29").Value = "Contains code examples for following concepts in TypeScript / ES6:`n1. Object literal short-hand notation in ES6`n2. Destructuring - Array and Object`n3. Modules, Import, Export"
$ws.Range("A30").Value = "store-app"
$ws.Range("B30").Value = "This app is created to understand:`n1. The organization of an Angular application`n2. How to install Bootstrap and how to include a reference to it within the app`n3. What is an Angular Module (ngModule)?`n4. What is a Component? What does a Component consist of?`n5. What is a Decorator? What is its purpose? Understand @NgModule and @Component decorators`n6. What is a root module? What is a root component?`n7. How to create components? How to compose simple components and create complex components`n8. What are component templates and styles? How to specify inline and external templates and styles for a component`n8. What are the various data binding techniques? How to use them?`n9. What is a Directive? What are the types? Examples of using built-in directives like ngModel, ngClass, ngStyle, ngFor, ngIf, etc.`n10. What is a Pipe? Examples of using built-in pipes like date, uppercase, currency. How to chain pipes?`nAt the end of day 1, participants will have a small, running Product Store application on their laptops."
$ws.Range("A31").Value = ""
$ws.Range("B31").Value = ""
$ws.Range("A32").Value = "day-03/"
$ws.Range("B32").Value = ""
$ws.Range("A33").Value = "lifecycle-hooks-app-start"
$ws.Range("B33").Value = "This app is used to understand and apply various Component Lifecyle Hooks (https://angular.io/guide/lifecycle-hooks). The start project provides a foundation to the participants to implement component lifecycle methods. It provides a hands-on experience to the participants who would like to code along with the instructor"
$ws.Range("A34").Value = "lifecycle-hooks-app-finish"
$ws.Range("B34").Value = "This app can be used as a reference by participants. It contains the code to demonstrate the usage of all lifecycle hooks."
$ws.Range("A35").Value = "store-app"
$ws.Range("B35").Value = "We continue with the day 2 'store-app' and understand & apply the following concepts:`n1. The need for creating small components`n2. How to refactor bigger components into smaller, more focussed components that follow 'Single Responsibility Principle'`n3. How do components interact with each other? (https://angular.io/guide/component-interaction)`n    a. Parent to child interaction`n    b. Child to parent interaction"
$ws.Range("A36").Value = "store-app-services"
$ws.Range("B36").Value = "We make a new copy the above app. This app will be used to understand:`n1. What is a Service? What is the need for a service?`n2. What is Dependency Injection (DI)? Why DI? What are the techniques to implement DI? `n    (https://angular.io/guide/dependency-injection)`n3. How to consume a service in a component?`n4. What is a provider?"
$ws.Range("A37").Value = ""
$ws.Range("B37").Value = ""
$ws.Range("A38").Value = "day-04/"
$ws.Range("B38").Value = ""
$ws.Range("A39").Value = "store-app-services"
$ws.Range("B39").Value = "We continue with day 3 'store-app-services' app to learn more about services:`n1. Understanding provider scope`n2. How a service can be used in component interaction?`nWe also learn:`n1. How to create a custom directive?`n2. How to create a custom pipe? "
$ws.Range("A40").Value = "store-app-forms-start"
$ws.Range("B40").Value = "This app is used to learn different types of form handling in Angular. It helps the participants to understand:`n1. Different approaches to handle user input through forms (https://angular.io/guide/forms-overview)`n    a. Template driven forms (https://angular.io/guide/forms)`n    b. Reactive forms (https://angular.io/guide/reactive-forms)`n2. The usage of various built-in validators`n3. How to create and use custom validators?`n4. How to create and use async validators?`n5. Various states of form and form controls`n6. How to set or retrieve values from form controls?`n7. Various directives used for template driven and reactive forms`nThe start project can be used by the participants to get hands-on coding experience along with the instructor."
$ws.Range("A41").Value = "store-app-forms-finish"
$ws.Range("B41").Value = "This app contains completed demo code for handling forms, both Template-driven and Reactive approaches. Participants can use this app as a reference"
$ws.Range("A42").Value = "store-app-routing-start"
$ws.Range("B42").Value = "This app is used to learn and apply routing and navigation within Angular apps. The Angular Router enables navigation from one view to another view (https://angular.io/guide/router). It helps the participant to:`n1. Define application routes and register them with the Router module`n2. Understand and use various routing related directives like <router-outlet>, routerLink, routerLinkActive, etc.`n3. Define wildcard route`n4. Define and retrieve route parameters`n5. Configure child routes`n6. Passing query parameters`n7. Understand the concept of route guards`n8. How to navigate programmatically using Router service?`nThe start project can be used by the participants to get hands-on coding experience along with the instructor."
$ws.Range("A43").Value = "store-app-routing-finish"
$ws.Range("B43").Value = "This app contains completed demo code for routing and navigation. Participants can use this app as a reference"
$ws.Range("A44").Value = ""
$ws.Range("B44").Value = ""
$ws.Range("A45").Value = "day-05/"
$ws.Range("B45").Value = ""
$ws.Range("A46").Value = "async-js"
$ws.Range("B46").Value = "This folder contains code samples for asynchronous JavaScript. The participant will learn different techniques to write async code in JavaScript:`n1. Using callbacks`n2. Using promises"
$ws.Range("A47").Value = "rxjs-examples"
$ws.Range("B47").Value = "This folder contains code samples for using RxJS Observables and various operators. The participant will learn:`n1. What is an Observable? Why Observable? (https://angular.io/guide/observables)`n2. Different techniques to create observables`n3. Apply numerous operators on observables"
$ws.Range("A48").Value = "store-app-http-start"
$ws.Range("B48").Value = "This app helps the participant to learn about server communication (HTTP) from Angular applications. This app helps participant to understand:`n1. How to communicate with the server via HTTP protocol using HttpClient service? (https://angular.io/guide/http)`n2. Getting JSON data from a server-side REST API`n3. Sending data to the server`n    a. Making a POST request`n    b. Making a DELETE request`n    c. Making a PUT/PATCH request`n    d. Adding headers`n4. Adding URL parameters`n5. How to handle errors?`n6. Understand the concept of interceptors`nThe start project can be used by the participants to get hands-on coding experience along with the instructor."
$ws.Range("A49").Value = "store-app-http-finish"
$ws.Range("B49").Value = "This app contains completed demo code for HTTP communcation with the server API. Participants can use this app as a reference"
$ws.Range("A50").Value = "store-app-http-server"
$ws.Range("B50").Value = "1. Contains server-api.json`n2. This file is used to simulate server-side REST API that is needed for the application`n3. 'json-server' tool is used to simulate a server. See '1-angular-workshop-lab-setup.pdf' document (docs folder) to know more about json-server"
$ws.Range("A51").Value = "store-app-unit-testing"
$ws.Range("B51").Value = "This app contains code for demonstrating various unit and integration test scenarios within an Angular app. (https://angular.io/guide/testing). The participant will learn about:`n1. Fundamentals of Angular testing`n2. Setup and Tear Down`n3. What are Spies? Why Spies?`n4. Code Coverage`n5. Angular Testing Utilities`n6. Testing Components`n7. Handling Component Dependencies`n8. Testing Async Operations"
$ws.Range("A52").Value = "store-app-unit-testing-server"
$ws.Range("B52").Value = "1. Contains server-api.json`n2. This file is used to simulate server-side REST API that is needed for the unit testing app`n3. 'json-server' tool is used to simulate a server. See '1-angular-workshop-lab-setup.pdf' document (docs folder) to know more about json-server"
$ws.Range("A53").Value = ""
$ws.Range("B53").Value = ""
$ws.Range("A54").Value = "misc-files"
$ws.Range("B54").Value = "Contains HTML markup for:`n1. Product Detail Component`n2. Product Form Component`nAlso, it contains Products JSON file used for simulating server side REST API"

# --- Step 3: row heights ---
$ws.Rows(25).RowHeight = 28.8
$ws.Rows(26).RowHeight = 57.6
$ws.Rows(29).RowHeight = 57.6
$ws.Rows(30).RowHeight = 201.6
$ws.Rows(33).RowHeight = 43.2
$ws.Rows(35).RowHeight = 86.4
$ws.Rows(36).RowHeight = 86.4
$ws.Rows(39).RowHeight = 100.8
$ws.Rows(40).RowHeight = 172.8
$ws.Rows(41).RowHeight = 28.8
$ws.Rows(42).RowHeight = 172.8
$ws.Rows(46).RowHeight = 57.6
$ws.Rows(47).RowHeight = 57.6
$ws.Rows(48).RowHeight = 201.6
$ws.Rows(49).RowHeight = 28.8
$ws.Rows(50).RowHeight = 57.6
$ws.Rows(51).RowHeight = 144
$ws.Rows(52).RowHeight = 57.6
$ws.Rows(54).RowHeight = 72

# --- Step 4: fix merged cell for the moved separator row (A42:B42 -> A44:B44) ---
$ws.Range("A42:B42").UnMerge() | Out-Null
$ws.Range("A44:B44").Merge() | Out-Null

# --- Step 5: update title-row selection to the merged title cell ---
$ws.Range("A1:B1").Select() | Out-Null
